$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Physiology")
$ws.Activate()
Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
